$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text format to cells whose new values look like plain numbers,
# so Excel keeps them as text (matching the source formatting) instead of
# converting them to floating point numbers.
$textCells = @('D5', 'D6', 'D8', 'D10', 'D12', 'D13', 'D18', 'D19', 'D21', 'D22', 'D23', 'D25', 'D26', 'D29', 'D31', 'D32', 'D33', 'D36', 'D40', 'D41')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '58.306.34'
$ws.Range('E2').Value = '  +2.24%  '
$ws.Range('D3').Value = '2.355.84'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '541.18'
$ws.Range('E5').Value = '  +1.55%  '
$ws.Range('D6').Value = '135.49'
$ws.Range('E6').Value = '  +2.42%  '
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('D8').Value = '0.562'
$ws.Range('E8').Value = '  +5.00%  '
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').Value = '5.60'
$ws.Range('E10').Value = '  +5.53%  '
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '0.354'
$ws.Range('E12').Value = '  +2.34%  '
$ws.Range('D13').Value = '23.86'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').Value = '2.775.31'
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '58.276.50'
$ws.Range('E15').Value = '  +2.11%  '
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '2.371.73'
$ws.Range('E17').Value = '  +2.28%  '
$ws.Range('D18').Value = '10.73'
$ws.Range('E18').Value = '  +3.03%  '
$ws.Range('D19').Value = '333.01'
$ws.Range('E19').Value = '  -1.26%  '
$ws.Range('E20').Value = '  +3.33%  '
$ws.Range('D21').Value = '6.80'
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').Value = '63.11'
$ws.Range('E23').Value = '  +2.40%  '
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '8.48'
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('E27').Value = '  +5.18%  '
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('D29').Value = '171.14'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = '0.0₃0737'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('D31').Value = '6.15'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('B32').Value = 'SuiNetwork'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D32').Value = '1.03'
$ws.Range('E32').Value = '  +13.51%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '18.46'
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  +7.26%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.78%  '
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('E38').Value = '  +4.27%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = '144.99'
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('D41').Value = '294.04'
$ws.Range('E41').Value = '  +4.66%  '
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('E43').Value = '  +1.32%  '
$ws.Range('E44').Value = '  +2.43%  '
$ws.Range('E45').Value = '  +2.41%  '
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('E51').Value = '  +0.53%  '
